{"js": "// Replace the lattice-multiplication exercise text in each table cell.\n// The table is a fixed 5-row x 3-column grid; every cell's text is\n// replaced in place (run formatting, e.g. <w:sz w:val=\"32\"/>, is kept\n// because we edit the existing paragraph range instead of clearing the\n// cell body).\nconst NEW_CELL_TEXT = [\n  \"32 x 66\\u000b  6    6\\u000b  ----\\u000b3|    |\\u000b2|    |\",\n  \"67 x 77\\u000b  7    7\\u000b  ----\\u000b6|    |\\u000b7|    |\",\n  \"62 x 85\\u000b  8    5\\u000b  ----\\u000b6|    |\\u000b2|    |\",\n  \"76 x 89\\u000b  8    9\\u000b  ----\\u000b7|    |\\u000b6|    |\",\n  \"51 x 35\\u000b  3    5\\u000b  ----\\u000b5|    |\\u000b1|    |\",\n  \"65 x 55\\u000b  5    5\\u000b  ----\\u000b6|    |\\u000b5|    |\",\n  \"26 x 86\\u000b  8    6\\u000b  ----\\u000b2|    |\\u000b6|    |\",\n  \"44 x 91\\u000b  9    1\\u000b  ----\\u000b4|    |\\u000b4|    |\",\n  \"57 x 44\\u000b  4    4\\u000b  ----\\u000b5|    |\\u000b7|    |\",\n  \"93 x 49\\u000b  4    9\\u000b  ----\\u000b9|    |\\u000b3|    |\",\n  \"73 x 23\\u000b  2    3\\u000b  ----\\u000b7|    |\\u000b3|    |\",\n  \"98 x 24\\u000b  2    4\\u000b  ----\\u000b9|    |\\u000b8|    |\",\n  \"70 x 30\\u000b  3    0\\u000b  ----\\u000b7|    |\\u000b0|    |\",\n  \"78 x 73\\u000b  7    3\\u000b  ----\\u000b7|    |\\u000b8|    |\",\n  \"82 x 51\\u000b  5    1\\u000b  ----\\u000b8|    |\\u000b2|    |\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\nlet index = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (index >= NEW_CELL_TEXT.length) {\n      break;\n    }\n    const cell = table.getCell(r, c);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const para = paragraphs.items[0];\n    const range = para.getRange();\n    range.insertText(NEW_CELL_TEXT[index], \"Replace\");\n    index++;\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Replace the lattice-multiplication exercise text in each table cell.\n# The table is a fixed 5-row x 3-column grid; every cell's text is\n# replaced in place via Range.Text so existing run formatting\n# (<w:sz w:val=\"32\"/>) on the cell's run is preserved.\n$newCellText = @(\n    \"32 x 66\" + [char]11 + \"  6    6\" + [char]11 + \"  ----\" + [char]11 + \"3|    |\" + [char]11 + \"2|    |\",\n    \"67 x 77\" + [char]11 + \"  7    7\" + [char]11 + \"  ----\" + [char]11 + \"6|    |\" + [char]11 + \"7|    |\",\n    \"62 x 85\" + [char]11 + \"  8    5\" + [char]11 + \"  ----\" + [char]11 + \"6|    |\" + [char]11 + \"2|    |\",\n    \"76 x 89\" + [char]11 + \"  8    9\" + [char]11 + \"  ----\" + [char]11 + \"7|    |\" + [char]11 + \"6|    |\",\n    \"51 x 35\" + [char]11 + \"  3    5\" + [char]11 + \"  ----\" + [char]11 + \"5|    |\" + [char]11 + \"1|    |\",\n    \"65 x 55\" + [char]11 + \"  5    5\" + [char]11 + \"  ----\" + [char]11 + \"6|    |\" + [char]11 + \"5|    |\",\n    \"26 x 86\" + [char]11 + \"  8    6\" + [char]11 + \"  ----\" + [char]11 + \"2|    |\" + [char]11 + \"6|    |\",\n    \"44 x 91\" + [char]11 + \"  9    1\" + [char]11 + \"  ----\" + [char]11 + \"4|    |\" + [char]11 + \"4|    |\",\n    \"57 x 44\" + [char]11 + \"  4    4\" + [char]11 + \"  ----\" + [char]11 + \"5|    |\" + [char]11 + \"7|    |\",\n    \"93 x 49\" + [char]11 + \"  4    9\" + [char]11 + \"  ----\" + [char]11 + \"9|    |\" + [char]11 + \"3|    |\",\n    \"73 x 23\" + [char]11 + \"  2    3\" + [char]11 + \"  ----\" + [char]11 + \"7|    |\" + [char]11 + \"3|    |\",\n    \"98 x 24\" + [char]11 + \"  2    4\" + [char]11 + \"  ----\" + [char]11 + \"9|    |\" + [char]11 + \"8|    |\",\n    \"70 x 30\" + [char]11 + \"  3    0\" + [char]11 + \"  ----\" + [char]11 + \"7|    |\" + [char]11 + \"0|    |\",\n    \"78 x 73\" + [char]11 + \"  7    3\" + [char]11 + \"  ----\" + [char]11 + \"7|    |\" + [char]11 + \"8|    |\",\n    \"82 x 51\" + [char]11 + \"  5    1\" + [char]11 + \"  ----\" + [char]11 + \"8|    |\" + [char]11 + \"2|    |\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$index = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($index -ge $newCellText.Length) {\n      break\n    }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newCellText[$index]\n    $index++\n  }\n}\n\n"}
